# edit.ps1
# Applies the edits described in the commit "Uploaded IRB approval letter"
# to the Research Information Sheet document.

$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    return $ok
}

# 1. Intro paragraph: extend the description of the study and add two line
#    breaks before "In the study, you will be asked ..."
$old1 = "a research study being performed by cognitive scientists in the Stanford Center for the Study of Language and Information. In this study, you will be asked questions about a variety of "
$new1 = "a research study being performed by cognitive scientists in the Stanford Center for the Study of Language and Information focused on how participants respond when they receive information about the study hypothesis. ^l^lIn the study, you will be asked questions about a variety of "
ReplaceText $old1 $new1

# 2. RISKS AND BENEFITS paragraph: "Your decision whether or not to
#    participate..." -- merge runs / clear grammar-check markers (no visible
#    text change).
$old2 = "Your decision whether or not to participate in this study will not "
$new2 = "Your decision whether or not to participate in this study will not "
ReplaceText $old2 $new2

# 3. PAYMENTS paragraph: change the compensation description.
$old3 = "You will receive partial course credit for participating in this study."
$new3 = "You will receive .50 SONA credits for participating in this study."
ReplaceText $old3 $new3

# 4. PARTICIPANT'S RIGHTS paragraph: "You may decline to answer any or all
#    of the following questions..." -- merge runs / clear grammar-check
#    markers (no visible text change).
$old4 = "You may decline to answer any or all of the following questions. You may decline further participation, at any time, without adverse consequences. "
$new4 = "You may decline to answer any or all of the following questions. You may decline further participation, at any time, without adverse consequences. "
ReplaceText $old4 $new4

# 5. CONTACT INFORMATION paragraph: "its procedures, risks and benefits,"
#    -- merge runs / clear grammar-check markers (no visible text change).
$old5 = "its procedures, risks and benefits, "
$new5 = "its procedures, risks and benefits, "
ReplaceText $old5 $new5

# 6. Footer: "File:TEM02C0" -- merge runs / clear grammar-check markers
#    (no visible text change).
$old6 = "File:TEM02C0"
$new6 = "File:TEM02C0"
ReplaceText $old6 $new6

# 7. Header text box: "Approval Date:   Monthname dd, 20yy Expiration
#    Date: Monthname dd, 20yy" -- merge runs / clear spell-check markers
#    (no visible text change).
$old7a = "Approval Date:   M"
$new7a = "Approval Date:   M"
ReplaceText $old7a $new7a

$old7b = "onthname dd, 20yy"
$new7b = "onthname dd, 20yy"
ReplaceText $old7b $new7b

$old7c = "Monthname dd, 20yy"
$new7c = "Monthname dd, 20yy"
ReplaceText $old7c $new7c
